$d = $word.ActiveDocument

# Update the header date line
$d.Content.Find.Execute("2026-01-31 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-01 Sunday", 2)

# Update the division problems in the table, addressed by (row, column)
# so that identical old values in different cells (e.g. "42÷6=" and "83÷9=")
# don't get cross-clobbered by a global text search/replace.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "42÷6=" },
    @{ Row = 1;  Col = 2; Text = "87÷8=" },
    @{ Row = 1;  Col = 3; Text = "13÷7=" },
    @{ Row = 1;  Col = 4; Text = "60÷9=" },
    @{ Row = 1;  Col = 5; Text = "17÷8=" },

    @{ Row = 5;  Col = 1; Text = "31÷2=" },
    @{ Row = 5;  Col = 2; Text = "29÷3=" },
    @{ Row = 5;  Col = 3; Text = "99÷8=" },
    @{ Row = 5;  Col = 4; Text = "54÷8=" },
    @{ Row = 5;  Col = 5; Text = "31÷6=" },

    @{ Row = 9;  Col = 1; Text = "92÷4=" },
    @{ Row = 9;  Col = 2; Text = "13÷9=" },
    @{ Row = 9;  Col = 3; Text = "83÷9=" },
    @{ Row = 9;  Col = 4; Text = "83÷3=" },
    @{ Row = 9;  Col = 5; Text = "93÷6=" },

    @{ Row = 13; Col = 1; Text = "90÷2=" },
    @{ Row = 13; Col = 2; Text = "12÷7=" },
    @{ Row = 13; Col = 3; Text = "16÷7=" },
    @{ Row = 13; Col = 4; Text = "24÷4=" },
    @{ Row = 13; Col = 5; Text = "61÷5=" },

    @{ Row = 17; Col = 1; Text = "51÷6=" },
    @{ Row = 17; Col = 2; Text = "45÷3=" },
    @{ Row = 17; Col = 3; Text = "89÷7=" },
    @{ Row = 17; Col = 4; Text = "75÷5=" },
    @{ Row = 17; Col = 5; Text = "49÷5=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
